# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stock) worksheet, populated for every existing data row, mirroring
# the legislator/report metadata that the crawler now attaches to each
# record in the exported dataframe.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "劉櫂豪"
$legislatorId = 1762
$reportDate = "2012-04-30"

# Find the last used data row (column A holds the running item number)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column H holds an ISO-looking date string ("2012-04-30"). Force it to
# text format up front so it is stored as a literal string instead of being
# reinterpreted as a date serial number; leave I (name, already text) and
# J (numeric legislator id) on the default General format.
$dateRange = $ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8))
$dateRange.NumberFormat = "@"

# Header row (row 1) — new columns H, I, J
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $reportDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
